$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 199
$ws.Range("B2").Value = "RG006719"
$ws.Range("C2").Value = "14YR785.7"
$ws.Range("E2").Value = "14K Rings- Dia Yellow XOXO ring"
$ws.Range("F2").Value = 7
$ws.Range("G2").Value = 1.685
$ws.Range("H2").Value = 0.118
$ws.Range("J2").Value = 1.661
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 348
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 586
